# Generate Report for Handoff
#
# Updates the localization-status report:
#   - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   - "Latest HO Xliff Generate Date" / "Latest Handback DateTime" timestamps bumped
#   - "Latest Handoff Datetime" timestamp (zh-cn) bumped
#   - Narrower "Status"/status columns on all three sheets

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Status changed from "Handed back: in sync with en-US" to "Ready for handoff"
#    on every sheet that shows it.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# 2. Refresh the handoff/handback timestamps.
$wsOverview.Range("G2").Value = "2016-08-25 09:04:05"
$wsDeDe.Range("H2").Value = "2016-08-25 09:04:05"
$wsZhCn.Range("H2").Value = "2016-08-25 09:03:56"

# 3. Narrow the status columns to fit the shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
